$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21767351.351937
$ws.Range("D2").Value = 71.15608

$ws.Range("B3").Value = 7952933.847096
$ws.Range("D3").Value = 12.998816
$ws.Range("E3").Value = 0.000004

$ws.Range("B4").Value = 101562096.096552
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = -182.121272
$ws.Range("H5").Value = -382.651956
$ws.Range("I5").Value = 18.409413
$ws.Range("J5").Value = 0.083815

$ws.Range("G6").Value = 159.274143
$ws.Range("H6").Value = -51.351063
$ws.Range("I6").Value = 369.899349
$ws.Range("J6").Value = 0.17778

$ws.Range("G7").Value = 341.395415
$ws.Range("H7").Value = 183.019328
$ws.Range("I7").Value = 499.771502
$ws.Range("J7").Value = 0.000002
